# feat: add 2022-Q4 data
#
# The existing "2022-Q3" sheet becomes "2022-Q4" (keeps its sheetId) and is
# populated with the new Q4 numbers; a fresh copy of the original "2022-Q3"
# sheet (with its original data untouched) is inserted right after it. The
# "总计" (totals) summary sheet gets its 2022-Q3 row's label/value corrected
# to 2022-Q4, plus a brand-new row re-adding the original 2022-Q3 totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write literal text into a cell without Excel's "looks like a
# number -> store as number" auto-coercion, and without touching any
# cell's number format / style (so no new style entries are minted).
# We do this by putting the text into a scratch cell via a formula that
# evaluates to a text literal, then Copy / PasteSpecial-values-only onto
# the destination (paste-values carries over the stored type as text but
# leaves the destination's existing formatting alone), then clean up the
# scratch cell.
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $ws = $range.Worksheet
    $scratch = $ws.Range("ZZ1000")
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.Clear() | Out-Null
}

# ---------------------------------------------------------------------
# 1. Duplicate the current "2022-Q3" sheet so the original data survives
#    on its own tab, then rename the two sheets: the original (which
#    keeps sheetId 2) becomes "2022-Q4", the new copy becomes "2022-Q3".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy([System.Reflection.Missing]::Value, $q3) | Out-Null

$q4sheet = $wb.Worksheets.Item("2022-Q3")        # original sheet (sheetId 2)
$q3copy  = $wb.Worksheets.Item("2022-Q3 (2)")    # fresh duplicate

$q4sheet.Name = "2022-Q4"
$q3copy.Name  = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Re-style the "2022-Q4" sheet to match the summary sheet's look
#    (header row + the index column use style index 2, same as "总计"),
#    reusing existing styles instead of creating new ones.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B1").Copy() | Out-Null
$q4sheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$total.Range("A2").Copy() | Out-Null
$q4sheet.Range("A2:A8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------
# 3. Overwrite "2022-Q4" with the new quarter's fund data.
# ---------------------------------------------------------------------
$q4sheet.Range("A2").Value = 0
$q4sheet.Range("A3").Value = 1
$q4sheet.Range("A4").Value = 2
$q4sheet.Range("A5").Value = 3
$q4sheet.Range("A6").Value = 4
$q4sheet.Range("A7").Value = 5
$q4sheet.Range("A8").Value = 6

Set-TextValue $q4sheet.Range("B2") "486001"
Set-TextValue $q4sheet.Range("C2") "工银瑞信中国机会全球配置股票（QDII）人民币"
Set-TextValue $q4sheet.Range("D2") "6.42"
Set-TextValue $q4sheet.Range("E2") "93.86"
Set-TextValue $q4sheet.Range("F2") "1.94"
Set-TextValue $q4sheet.Range("G2") "0.1245"
$q4sheet.Range("H2").Value = 6

Set-TextValue $q4sheet.Range("B3") "009562"
Set-TextValue $q4sheet.Range("C3") "工银全球股票（QDII）美元"
Set-TextValue $q4sheet.Range("D3") "6.42"
Set-TextValue $q4sheet.Range("E3") "93.86"
Set-TextValue $q4sheet.Range("F3") "1.94"
Set-TextValue $q4sheet.Range("G3") "0.1245"
$q4sheet.Range("H3").Value = 6

Set-TextValue $q4sheet.Range("B4") "009563"
Set-TextValue $q4sheet.Range("C4") "工银全球股票（QDII）港币"
Set-TextValue $q4sheet.Range("D4") "6.42"
Set-TextValue $q4sheet.Range("E4") "93.86"
Set-TextValue $q4sheet.Range("F4") "1.94"
Set-TextValue $q4sheet.Range("G4") "0.1245"
$q4sheet.Range("H4").Value = 6

Set-TextValue $q4sheet.Range("B5") "012751"
Set-TextValue $q4sheet.Range("C5") "建信纳斯达克100指数（QDII）A 美元现汇"
Set-TextValue $q4sheet.Range("D5") "1.06"
Set-TextValue $q4sheet.Range("E5") "82.28"
Set-TextValue $q4sheet.Range("F5") "11.61"
Set-TextValue $q4sheet.Range("G5") "0.1231"
$q4sheet.Range("H5").Value = 1

Set-TextValue $q4sheet.Range("B6") "012752"
Set-TextValue $q4sheet.Range("C6") "建信纳斯达克100指数（QDII）C 人民币"
Set-TextValue $q4sheet.Range("D6") "1.06"
Set-TextValue $q4sheet.Range("E6") "82.28"
Set-TextValue $q4sheet.Range("F6") "11.61"
Set-TextValue $q4sheet.Range("G6") "0.1231"
$q4sheet.Range("H6").Value = 1

Set-TextValue $q4sheet.Range("B7") "012753"
Set-TextValue $q4sheet.Range("C7") "建信纳斯达克100指数（QDII）C 美元现汇"
Set-TextValue $q4sheet.Range("D7") "1.06"
Set-TextValue $q4sheet.Range("E7") "82.28"
Set-TextValue $q4sheet.Range("F7") "11.61"
Set-TextValue $q4sheet.Range("G7") "0.1231"
$q4sheet.Range("H7").Value = 1

Set-TextValue $q4sheet.Range("B8") "486002"
Set-TextValue $q4sheet.Range("C8") "工银全球精选股票（QDII）"
Set-TextValue $q4sheet.Range("D8") "3.92"
Set-TextValue $q4sheet.Range("E8") "94.38"
Set-TextValue $q4sheet.Range("F8") "3.06"
Set-TextValue $q4sheet.Range("G8") "0.1200"
$q4sheet.Range("H8").Value = 1

# ---------------------------------------------------------------------
# 4. Update the "总计" summary sheet: the old 2022-Q3 row now reports
#    2022-Q4 numbers, and a new row re-adds the original 2022-Q3 totals.
# ---------------------------------------------------------------------
Set-TextValue $total.Range("B2") "2022-Q4"
$total.Range("D2").Value = 0.86

$total.Range("A2").Copy() | Out-Null
$total.Range("A3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$total.Range("A3").Value = 1
Set-TextValue $total.Range("B3") "2022-Q3"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 0.87

# ---------------------------------------------------------------------
# 5. Leave the active tab on "总计", matching the original workbook.
# ---------------------------------------------------------------------
$total.Activate() | Out-Null
